$d = $word.ActiveDocument

$replacements = @(
    @("36÷4=", "37÷3="),
    @("15÷7=", "30÷9="),
    @("86÷7=", "92÷7="),
    @("81÷5=", "82÷7="),
    @("40÷9=", "81÷6="),
    @("81÷8=", "99÷8="),
    @("11÷2=", "51÷6="),
    @("32÷4=", "34÷7="),
    @("12÷9=", "66÷7="),
    @("97÷8=", "65÷8="),
    @("14÷4=", "34÷6="),
    @("66÷2=", "45÷8="),
    @("15÷5=", "85÷4="),
    @("64÷6=", "33÷8="),
    @("14÷9=", "29÷9="),
    @("97÷5=", "70÷5="),
    @("39÷9=", "35÷6="),
    @("42÷6=", "31÷5="),
    @("22÷3=", "63÷4="),
    @("50÷8=", "13÷6="),
    @("88÷4=", "95÷8="),
    @("59÷6=", "83÷2="),
    @("51÷5=", "85÷6="),
    @("88÷9=", "44÷9="),
    @("65÷2=", "96÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
